# "create new property completed"
# Adds three new reference sheets (constructionTypeDetails, floorDetails,
# approvalDetails) ahead of propertyHeaderDetails, fills them with data,
# nudges a date + selection on assessmentDetails, and moves the active tab
# from amenities to the new floorDetails sheet.

$wb = $excel.ActiveWorkbook

# Helper: Excel's xlsx COLUMN width attribute is ColumnWidth + 0.8333333333333334
# (measured empirically against this runtime's AutoFit/width pipeline), so we
# back-solve from the target xlsx "width" values recorded in the authored file.
function Set-ColWidth($range, [double]$xmlWidth) {
    $range.ColumnWidth = $xmlWidth - 0.8333333333333334
}

# ---------------------------------------------------------------------------
# 1. assessmentDetails: bump the registration date and move the selection
# ---------------------------------------------------------------------------
$assessment = $wb.Worksheets.Item("assessmentDetails")
$assessment.Range("F2").Value = 42350
$assessment.Activate()
$assessment.Range("F2").Select()

# ---------------------------------------------------------------------------
# 2. Insert the three new sheets right after "amenities" (and therefore right
#    before "propertyHeaderDetails", which gets pushed to the end).
# ---------------------------------------------------------------------------
$amenities = $wb.Worksheets.Item("amenities")

$constructionTypeDetails = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $amenities)
$constructionTypeDetails.Name = "constructionTypeDetails"

$floorDetails = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $constructionTypeDetails)
$floorDetails.Name = "floorDetails"

$approvalDetails = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $floorDetails)
$approvalDetails.Name = "approvalDetails"

# ---------------------------------------------------------------------------
# 3. constructionTypeDetails content
#    (cell fill order matches the authored shared-strings sequence)
# ---------------------------------------------------------------------------
Set-ColWidth $constructionTypeDetails.Columns("A:A") 21.33203125

$constructionTypeDetails.Range("A1").Value = "dataName"
$constructionTypeDetails.Range("B1").Value = "floorType"
$constructionTypeDetails.Range("C1").Value = "roofType"
$constructionTypeDetails.Range("D1").Value = "wallType"
$constructionTypeDetails.Range("E1").Value = "woodType"

$constructionTypeDetails.Range("A2").Value = "defaultConstructionType"
$constructionTypeDetails.Range("B2").Value = "Black Stones"
$constructionTypeDetails.Range("C2").Value = "Absheet"
$constructionTypeDetails.Range("D2").Value = "BAMBOO"
$constructionTypeDetails.Range("E2").Value = "Allmixing"

$constructionTypeDetails.Range("E3").Select()

# ---------------------------------------------------------------------------
# 4. floorDetails content
# ---------------------------------------------------------------------------
Set-ColWidth $floorDetails.Columns("B:B") 11.6640625
Set-ColWidth $floorDetails.Columns("C:C") 20
Set-ColWidth $floorDetails.Columns("D:D") 13.33203125
Set-ColWidth $floorDetails.Columns("G:G") 13.5
Set-ColWidth $floorDetails.Columns("H:H") 15.1640625
Set-ColWidth $floorDetails.Columns("I:I") 16.33203125
Set-ColWidth $floorDetails.Columns("J:J") 15.6640625
Set-ColWidth $floorDetails.Columns("K:K") 6.33203125
Set-ColWidth $floorDetails.Columns("L:L") 7.5
Set-ColWidth $floorDetails.Columns("M:M") 23.1640625
Set-ColWidth $floorDetails.Columns("N:N") 20.33203125
Set-ColWidth $floorDetails.Columns("O:O") 21.1640625

# NOTE: the fill order below intentionally mirrors the authored workbook's
# shared-strings allocation sequence (e.g. C1 before B1, and A2/C2/B2 before
# D2/E2/F2) rather than simple left-to-right order.
$floorDetails.Range("A1").Value = "dataName"
$floorDetails.Range("G2").Value = "Bimal"
$floorDetails.Range("C1").Value = "classificationOfBuilding"
$floorDetails.Range("B1").Value = "floorNumber"
$floorDetails.Range("D1").Value = "natureOfUsage"
$floorDetails.Range("E1").Value = "firmName"
$floorDetails.Range("F1").Value = "occupancy"
$floorDetails.Range("G1").Value = "occupantName"
$floorDetails.Range("H1").Value = "constructionDate"

$floorDetails.Range("A2").Value = "firstFloor"
$floorDetails.Range("C2").Value = "Huts"
$floorDetails.Range("B2").Value = "1st floor"
$floorDetails.Range("D2").Value = "Residence"
$floorDetails.Range("E2").Value = "NA"
$floorDetails.Range("F2").Value = "Owner"

$floorDetails.Range("I1").Value = "effectiveFromDate"
$floorDetails.Range("J1").Value = "unstructuredLand"
$floorDetails.Range("K1").Value = "length"
$floorDetails.Range("L1").Value = "breadth"
$floorDetails.Range("M1").Value = "buildingPermissionNumber"
$floorDetails.Range("N1").Value = "buildingPermissionDate"
$floorDetails.Range("O1").Value = "plinthAreaInBuildingPlan"

$floorDetails.Range("J2").Value = "No"
$floorDetails.Range("M2").Value = "11/22"

$floorDetails.Range("K2").Value = 10
$floorDetails.Range("L2").Value = 20
$floorDetails.Range("O2").Value = 30

# Dates: copy the existing short-date format (already used on
# assessmentDetails!F2) onto the new date cells, then stamp in the values.
$floorDetails.Range("H2").Value = 42653
$floorDetails.Range("I2").Value = 42654
$floorDetails.Range("N2").Value = 42358
$assessment.Range("F2").Copy()
$floorDetails.Range("H2,I2,N2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$floorDetails.Range("N2").Select()
$floorDetails.Activate()

# ---------------------------------------------------------------------------
# 5. approvalDetails content
# ---------------------------------------------------------------------------
Set-ColWidth $approvalDetails.Columns("A:A") 14.5
Set-ColWidth $approvalDetails.Columns("B:B") 18.33203125
Set-ColWidth $approvalDetails.Columns("C:C") 18.1640625
Set-ColWidth $approvalDetails.Columns("D:D") 27
Set-ColWidth $approvalDetails.Columns("E:E") 15.5

# NOTE: fill order again mirrors the authored shared-strings sequence: the
# data row's B2:D2 values were entered before the A2 row-key and E2 remark.
$approvalDetails.Range("A1").Value = "dataName"
$approvalDetails.Range("B1").Value = "approverDepartment"
$approvalDetails.Range("C1").Value = "approverDesignation"
$approvalDetails.Range("D1").Value = "approver"
$approvalDetails.Range("E1").Value = "approverRemarks"

$approvalDetails.Range("B2").Value = "REVENUE"
$approvalDetails.Range("C2").Value = "Bill Collector"
$approvalDetails.Range("D2").Value = "D.Khasim ~ REV_Bill Collector_1"
$approvalDetails.Range("A2").Value = "defaultApprover"
$approvalDetails.Range("E2").Value = "Forward to BC"

$approvalDetails.Range("E3").Select()

# ---------------------------------------------------------------------------
# 6. Final active tab is floorDetails (matches workbookView activeTab="6").
# ---------------------------------------------------------------------------
$floorDetails.Activate()
$floorDetails.Range("N2").Select()
